# Fix bug of WireBuilder lock:
#  - uart_rx sheet: Width of m_axis_tdata port (C5) was 0, should be 1
#  - uart_tx sheet: Width of s_axis_tdata port (C5) was 0, should be 1
# Also restore the author's last on-screen selection/active sheet state.

$wb = $excel.ActiveWorkbook

$wsRx = $wb.Worksheets.Item("uart_rx")
$wsTx = $wb.Worksheets.Item("uart_tx")
$wsUart = $wb.Worksheets.Item("uart")

$wsRx.Range("C5").Value = 1
$wsTx.Range("C5").Value = 1

# Leave "uart" sheet's selection at A1
$wsUart.Activate()
$wsUart.Range("A1").Select()

# uart_tx sheet selection left on C5 (the cell that was just edited)
$wsTx.Activate()
$wsTx.Range("C5").Select()

# uart_rx ends up as the active sheet/tab, with D25 selected
$wsRx.Activate()
$wsRx.Range("D25").Select()
